# F030. Split PrEP Given and PrEP Education in Services received export
#
# The header row (row 1) of the Services_Received sheet had a single
# "PrEP" column. It needs to become two columns:
#   - the existing column is renamed to "PrEP Given"
#   - a new column "PrEP Education" is inserted immediately after it
#
# Everything to the right of "PrEP" shifts one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the existing "PrEP" header cell in row 1.
$prepCell = $ws.Rows.Item(1).Find("PrEP")

# Insert a new column right after it, shifting the remaining columns right.
$nextCol = $prepCell.Offset(0, 1).EntireColumn
$nextCol.Insert()

# Rename the original column's header and set the new column's header.
$prepCell.Value = "PrEP Given"
$prepCell.Offset(0, 1).Value = "PrEP Education"
